$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dRng = $ws.Range("D2")
$dRng.NumberFormat = "@"
$dRng.Value = "27.812.05"
$dRng.Style = "Normal"
$ws.Range("E2").Value = "  +0.70%  "

$dRng = $ws.Range("D3")
$dRng.NumberFormat = "@"
$dRng.Value = "1.858.51"
$dRng.Style = "Normal"
$ws.Range("E3").Value = "  +0.39%  "

$dRng = $ws.Range("D4")
$dRng.NumberFormat = "@"
$dRng.Value = "1.037"
$dRng.Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "

$dRng = $ws.Range("D5")
$dRng.NumberFormat = "@"
$dRng.Value = "323.42"
$dRng.Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "

$dRng = $ws.Range("D6")
$dRng.NumberFormat = "@"
$dRng.Value = "1.033"
$dRng.Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

$dRng = $ws.Range("D7")
$dRng.NumberFormat = "@"
$dRng.Value = "0.4416"
$dRng.Style = "Normal"
$ws.Range("E7").Value = "  +0.83%  "

$dRng = $ws.Range("D8")
$dRng.NumberFormat = "@"
$dRng.Value = "0.3823"
$dRng.Style = "Normal"
$ws.Range("E8").Value = "  +1.85%  "

$dRng = $ws.Range("D9")
$dRng.NumberFormat = "@"
$dRng.Value = "0.07444"
$dRng.Style = "Normal"
$ws.Range("E9").Value = "  +0.51%  "

$dRng = $ws.Range("D10")
$dRng.NumberFormat = "@"
$dRng.Value = "0.8885"
$dRng.Style = "Normal"
$ws.Range("E10").Value = "  +1.50%  "

$dRng = $ws.Range("D11")
$dRng.NumberFormat = "@"
$dRng.Value = "21.65"
$dRng.Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "

$dRng = $ws.Range("D12")
$dRng.NumberFormat = "@"
$dRng.Value = "1.873.12"
$dRng.Style = "Normal"
$ws.Range("E12").Value = "  +0.59%  "

$dRng = $ws.Range("D13")
$dRng.NumberFormat = "@"
$dRng.Value = "5.537"
$dRng.Style = "Normal"
$ws.Range("E13").Value = "  +0.31%  "

$dRng = $ws.Range("D14")
$dRng.NumberFormat = "@"
$dRng.Value = "6.751"
$dRng.Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "

$dRng = $ws.Range("D15")
$dRng.NumberFormat = "@"
$dRng.Value = "0.07212"
$dRng.Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

$dRng = $ws.Range("D16")
$dRng.NumberFormat = "@"
$dRng.Value = "86.03"
$dRng.Style = "Normal"
$ws.Range("E16").Value = "  +3.75%  "

$dRng = $ws.Range("D17")
$dRng.NumberFormat = "@"
$dRng.Value = "1.039"
$dRng.Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "

$dRng = $ws.Range("D18")
$dRng.NumberFormat = "@"
$dRng.Value = "0.000009108"
$dRng.Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("E19").Value = "  +0.34%  "

$dRng = $ws.Range("D20")
$dRng.NumberFormat = "@"
$dRng.Value = "15.59"
$dRng.Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "

$dRng = $ws.Range("D21")
$dRng.NumberFormat = "@"
$dRng.Value = "27.828.33"
$dRng.Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "

$dRng = $ws.Range("D22")
$dRng.NumberFormat = "@"
$dRng.Value = "5.302"
$dRng.Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("E23").Value = "  +0.53%  "

$dRng = $ws.Range("D24")
$dRng.NumberFormat = "@"
$dRng.Value = "2.095.70"
$dRng.Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "

$dRng = $ws.Range("D25")
$dRng.NumberFormat = "@"
$dRng.Value = "2.066"
$dRng.Style = "Normal"
$ws.Range("E25").Value = "  +6.53%  "

$dRng = $ws.Range("D26")
$dRng.NumberFormat = "@"
$dRng.Value = "159.15"
$dRng.Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "

$dRng = $ws.Range("D27")
$dRng.NumberFormat = "@"
$dRng.Value = "18.79"
$dRng.Style = "Normal"
$ws.Range("E27").Value = "  +0.19%  "

$dRng = $ws.Range("D28")
$dRng.NumberFormat = "@"
$dRng.Value = "2.013"
$dRng.Style = "Normal"
$ws.Range("E28").Value = "  +3.95%  "

$dRng = $ws.Range("D29")
$dRng.NumberFormat = "@"
$dRng.Value = "5.372"
$dRng.Style = "Normal"
$ws.Range("E29").Value = "  +1.51%  "

$dRng = $ws.Range("D30")
$dRng.NumberFormat = "@"
$dRng.Value = "118.64"
$dRng.Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "

$dRng = $ws.Range("D31")
$dRng.NumberFormat = "@"
$dRng.Value = "0.09114"
$dRng.Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "

$dRng = $ws.Range("D32")
$dRng.NumberFormat = "@"
$dRng.Value = "0.7758"
$dRng.Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "

$ws.Range("E33").Value = "  +0.74%  "

$dRng = $ws.Range("D34")
$dRng.NumberFormat = "@"
$dRng.Value = "3.025"
$dRng.Style = "Normal"
$ws.Range("E34").Value = "  +4.84%  "

$dRng = $ws.Range("D35")
$dRng.NumberFormat = "@"
$dRng.Value = "4.616"
$dRng.Style = "Normal"
$ws.Range("E35").Value = "  +2.31%  "

$dRng = $ws.Range("D36")
$dRng.NumberFormat = "@"
$dRng.Value = "1.035"
$dRng.Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "

$dRng = $ws.Range("D37")
$dRng.NumberFormat = "@"
$dRng.Value = "1.157"
$dRng.Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

$dRng = $ws.Range("D38")
$dRng.NumberFormat = "@"
$dRng.Value = "0.01986"
$dRng.Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "

$dRng = $ws.Range("D39")
$dRng.NumberFormat = "@"
$dRng.Value = "0.05331"
$dRng.Style = "Normal"
$ws.Range("E39").Value = "  +0.96%  "

$dRng = $ws.Range("D40")
$dRng.NumberFormat = "@"
$dRng.Value = "2.875"
$dRng.Style = "Normal"
$ws.Range("E40").Value = "  +2.18%  "

$dRng = $ws.Range("D41")
$dRng.NumberFormat = "@"
$dRng.Value = "0.5222"
$dRng.Style = "Normal"
$ws.Range("E41").Value = "  +0.84%  "

$dRng = $ws.Range("D42")
$dRng.NumberFormat = "@"
$dRng.Value = "6.945"
$dRng.Style = "Normal"
$ws.Range("E42").Value = "  +3.18%  "

$dRng = $ws.Range("D43")
$dRng.NumberFormat = "@"
$dRng.Value = "0.1680"
$dRng.Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "

$dRng = $ws.Range("D44")
$dRng.NumberFormat = "@"
$dRng.Value = "8.810"
$dRng.Style = "Normal"
$ws.Range("E44").Value = "  +2.73%  "

$dRng = $ws.Range("D45")
$dRng.NumberFormat = "@"
$dRng.Value = "111.08"
$dRng.Style = "Normal"
$ws.Range("E45").Value = "  +1.81%  "

$dRng = $ws.Range("D46")
$dRng.NumberFormat = "@"
$dRng.Value = "10.79"
$dRng.Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "

$ws.Range("E47").Value = "  +0.39%  "

$dRng = $ws.Range("D48")
$dRng.NumberFormat = "@"
$dRng.Value = "0.06589"
$dRng.Style = "Normal"
$ws.Range("E48").Value = "  +2.95%  "

$ws.Range("E49").Value = "  -0.04%  "

$dRng = $ws.Range("D50")
$dRng.NumberFormat = "@"
$dRng.Value = "0.4738"
$dRng.Style = "Normal"
$ws.Range("E50").Value = "  +1.69%  "

$dRng = $ws.Range("D51")
$dRng.NumberFormat = "@"
$dRng.Value = "1.896"
$dRng.Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
